# Hoan thanh them tab QuanLyDeThi, frmLogin
# (workbook-visible portion of that commit: add maMH/maKhoi columns to the
# CauHoi question table and extend Table1 to cover them)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CauHoi")
$ws.Activate()

$lastRow = 21

# --- add the "maMH" column in D, filled with "T" for every question row ---
$lo = $ws.ListObjects.Item(1)
$colMH = $lo.ListColumns.Add()
$ws.Range("D1").Value = "maMH"
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 4).Value = "T"
}

# --- add the "maKhoi" column in E, filled with "K10" for every question row ---
$colKhoi = $lo.ListColumns.Add()
$ws.Range("E1").Value = "maKhoi"
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 5).Value = "K10"
}

# make sure the table now spans A1:E21 with both new columns included
$lo.Resize($ws.Range("A1:E21"))

# match the author's final view state: scrolled right with the new
# maKhoi column selected
$ws.Range("E2:E21").Select()
